$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "cards"
$ws.Range("E2").Value = "Permament Residence Permit"
$ws.Range("L2").Value = 15
$ws.Range("N2").Value = 5

$ws.Range("O2").Select()
